$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record row at position 181 (pushes the old 181..238
# rows down to 182..239, and the dimension grows from R238 to R239).
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row with the new data record.
$ws.Range("A181").Value = 8
$ws.Range("B181").Value = "Terminal La Palmera de La Serena"
$ws.Range("C181").Value = "Coquimbo"
$ws.Range("D181").Value = Get-Date -Year 2022 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Range("E181").Value = 4
$ws.Range("F181").Value = 100112031
$ws.Range("G181").Value = "Poroto verde"
$ws.Range("H181").Value = "Magnum"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 400
$ws.Range("K181").Value = 25000
$ws.Range("L181").Value = 26000
$ws.Range("M181").Value = 25500
$ws.Range("N181").Value = "`$/malla 25 kilos"
$ws.Range("O181").Value = "Perú"
$ws.Range("P181").Value = 1020
$ws.Range("Q181").Value = 25
$ws.Range("R181").Value = "Hortaliza"
